$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.186.34"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.563.72"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.09"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.70"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.63"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.65"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.022.98"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.113.67"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  +4.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.563.93"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.40"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.44"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.81"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.69"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.65"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.684.04"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.09"
$ws.Range("E27").Value = "  +11.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.57"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.00"
$ws.Range("E31").Value = "  +8.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0827"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "443.54"
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.86"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.34"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.52"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "151.39"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.16"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0554"
$ws.Range("E45").Value = "  +6.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.615"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0977"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0242"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.45"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("E51").Value = "  -0.25%  "
